# Daily attendance processing - 2025-10-17 07:19:45
# For the rows below, column G holds a two-item, comma-separated list of
# "last edited by" actors (e.g. "System, someone@example.com"). This swaps
# the order of the two items for each listed row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(3,4,6,7,10,11,12,13,14,15,17,18,19,30,31,33,34,37,38,39,40,41,42,44,45,46,57,58,60,61,64,65,66,67,68,69,71,72,73,83,86,87,88,89,90,93,95,96,97,99,109,112,113,114,115,116,119,121,122,123,125,135,138,139,140,141,142,145,147,148,149,151)

foreach ($row in $rows) {
    $cell = $ws.Cells.Item($row, 7)   # column G
    $current = [string]$cell.Value2
    $parts = $current.Split(",")
    if ($parts.Count -eq 2) {
        $first = $parts[0].Trim()
        $second = $parts[1].Trim()
        $cell.Value = "$second, $first"
    }
}
